$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value2 = $ws.Range("B1").Value2
$ws.Range("C2").Value2 = $ws.Range("B2").Value2
$ws.Range("C3").Value2 = $ws.Range("B3").Value2
$ws.Range("C4").Value2 = $ws.Range("B4").Value2
$ws.Range("C5").Value2 = $ws.Range("B5").Value2
$ws.Range("C7").Value2 = $ws.Range("B7").Value2

$ws.Range("L17").Select()
